$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 7 (new "date_example" / "prompt_example" demo blocks).
# Excel copies the row-6 formatting down into the freshly inserted rows, which already
# matches the target column styles for columns A-G.
$ws.Rows("7:8").Insert()

# Row 5: the "age" answer now routes to the new date-picker demo instead of "self".
$ws.Range("E5").Value = "GO(date_example)"
$ws.Rows("5").RowHeight = 13.95

# New row 7: date_example question (type D = date picker demo)
$ws.Range("A7").Value = "date_example"
$ws.Range("B7").Value = "Date example"
$ws.Range("C7").Value = "D"
$ws.Range("D7").Value = "Give me your birthdate. I promise to keep the secret :)"
$ws.Range("E7").Value = "GO(prompt_example)"
$ws.Range("H7").Value = "What is this?"
$ws.Range("I7").Value = "An example of date picker component"
$ws.Rows("7").RowHeight = 13.95

# New row 8: prompt_example question (type F = free text demo)
$ws.Range("A8").Value = "prompt_example"
$ws.Range("B8").Value = "Prompt example"
$ws.Range("C8").Value = "F"
$ws.Range("D8").Value = "Now tell me what is you opinion on medical chatbots!"
$ws.Range("E8").Value = "GO(self)"
$ws.Range("H8").Value = "What is this?"
$ws.Range("I8").Value = "An example of date picker component"
$ws.Rows("8").RowHeight = 14.9

# Row 9 (previously row 7 - the "self" question) gains a Hint/Popup pair and its row
# grows slightly taller to fit the new popup text.
$ws.Range("H9").Value = "What is this?"
$ws.Range("I9").Value = "An example of free text input component"
$ws.Rows("9").RowHeight = 28.9

# Column widths: new column A, and a new column I sized for the popup text.
$ws.Columns("A").ColumnWidth = 14
$ws.Columns("I").ColumnWidth = 30.35

# View: scroll back to the top-left and select E9 (the new "self" query output cell).
[void]$ws.Range("E9").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
